$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Saldo_guide"

# Bulk update Dt. Referencia (column G) for all data rows: 2024-03-27 -> 2024-04-01 (serial 45378 -> 45383)
$ws.Range("G2:G310").Value2 = 45383

# Refreshed balances for accounts with new transaction activity
$ws.Range("D2").Value2 = 266.62
$ws.Range("H2").Value2 = 266.62

$ws.Range("D5").Value2 = 59645.52
$ws.Range("E5").Value2 = -59721.43
$ws.Range("H5").Value2 = -75.91

$ws.Range("D6").Value2 = 127459.74
$ws.Range("H6").Value2 = 127459.74

$ws.Range("D8").Value2 = 0.04
$ws.Range("H8").Value2 = 0.04

$ws.Range("D12").Value2 = 105801.91
$ws.Range("H12").Value2 = 105801.91

$ws.Range("D23").Value2 = 104934.56
$ws.Range("E23").Value2 = -113482.42
$ws.Range("H23").Value2 = -8547.86

$ws.Range("D27").Value2 = 53672.1
$ws.Range("E27").Value2 = -6419.1
$ws.Range("H27").Value2 = 47253

$ws.Range("D35").Value2 = 214.17
$ws.Range("H35").Value2 = 214.17

$ws.Range("D47").Value2 = 144530.26999999999
$ws.Range("H47").Value2 = 144530.26999999999

$ws.Range("D55").Value2 = 31140.15
$ws.Range("E55").Value2 = 0
$ws.Range("H55").Value2 = 31140.15

$ws.Range("D61").Value2 = 49195.25
$ws.Range("E61").Value2 = -7802.93
$ws.Range("H61").Value2 = 41392.32

$ws.Range("D63").Value2 = 825.44
$ws.Range("H63").Value2 = 825.44

$ws.Range("D64").Value2 = 1.54
$ws.Range("H64").Value2 = 1.54

$ws.Range("D67").Value2 = 14.37
$ws.Range("H67").Value2 = 14.37

$ws.Range("D69").Value2 = 121261.09
$ws.Range("H69").Value2 = 121261.09

$ws.Range("D71").Value2 = 879.57
$ws.Range("H71").Value2 = 879.57

$ws.Range("D73").Value2 = 22.95
$ws.Range("H73").Value2 = 22.95

$ws.Range("D76").Value2 = 8079.07
$ws.Range("E76").Value2 = -8123.57
$ws.Range("H76").Value2 = -44.5

$ws.Range("D78").Value2 = 7142.86
$ws.Range("E78").Value2 = 5535.28
$ws.Range("H78").Value2 = 12678.14

$ws.Range("D101").Value2 = 14.51
$ws.Range("H101").Value2 = 14.51

$ws.Range("D111").Value2 = 64.84
$ws.Range("H111").Value2 = 64.84

$ws.Range("D115").Value2 = 34834.980000000003
$ws.Range("E115").Value2 = 0
$ws.Range("H115").Value2 = 34834.980000000003

$ws.Range("D117").Value2 = 893.38
$ws.Range("H117").Value2 = 893.38

$ws.Range("D118").Value2 = 10000.09
$ws.Range("H118").Value2 = 10000.09

$ws.Range("D119").Value2 = 1.99
$ws.Range("H119").Value2 = 1.99

$ws.Range("D120").Value2 = 83564.66
$ws.Range("E120").Value2 = 0
$ws.Range("H120").Value2 = 83564.66

$ws.Range("D121").Value2 = 405.84
$ws.Range("H121").Value2 = 405.84

$ws.Range("D124").Value2 = 10739.86
$ws.Range("H124").Value2 = 10739.86

$ws.Range("D125").Value2 = 97.8
$ws.Range("H125").Value2 = 97.8

$ws.Range("D126").Value2 = 875.78
$ws.Range("H126").Value2 = 875.78

$ws.Range("D129").Value2 = 0.02
$ws.Range("H129").Value2 = 0.02

$ws.Range("D136").Value2 = 0
$ws.Range("H136").Value2 = 0

$ws.Range("D137").Value2 = 3488.21
$ws.Range("H137").Value2 = 3488.21

$ws.Range("D151").Value2 = 41896.379999999997
$ws.Range("E151").Value2 = -4796.8599999999997
$ws.Range("H151").Value2 = 37099.519999999997

$ws.Range("D160").Value2 = 31077.35
$ws.Range("H160").Value2 = 31077.35

$ws.Range("D167").Value2 = 54425.81
$ws.Range("H167").Value2 = 54425.81

$ws.Range("D175").Value2 = 812.28
$ws.Range("H175").Value2 = 812.28

$ws.Range("D184").Value2 = 10911.45
$ws.Range("H184").Value2 = 10911.45

$ws.Range("D187").Value2 = 230.91
$ws.Range("H187").Value2 = 230.91

$ws.Range("D188").Value2 = 0.22
$ws.Range("H188").Value2 = 0.22

$ws.Range("D192").Value2 = 31587.84
$ws.Range("H192").Value2 = 31587.84

$ws.Range("D200").Value2 = 10725.59
$ws.Range("E200").Value2 = -2027.59
$ws.Range("H200").Value2 = 8698

$ws.Range("D201").Value2 = 1.9
$ws.Range("H201").Value2 = 1.9

$ws.Range("D251").Value2 = 12.63
$ws.Range("H251").Value2 = 12.63

$ws.Range("D267").Value2 = 3157.16
$ws.Range("E267").Value2 = -3132.5
$ws.Range("H267").Value2 = 24.66

$ws.Range("D268").Value2 = 11.45
$ws.Range("H268").Value2 = 11.45

$ws.Range("D274").Value2 = 9111.83
$ws.Range("E274").Value2 = -9128.25
$ws.Range("H274").Value2 = -16.420000000000002

$ws.Range("D290").Value2 = 466.72
$ws.Range("H290").Value2 = 466.72
